$wb = $excel.ActiveWorkbook

# Template sheet to clone the formatting/layout from for each new market sheet.
$template = $wb.Worksheets.Item("Denmark")

function Add-MarketSheet {
    param(
        [string]$Name,
        [string]$TicketRef,
        [string]$MarketLabel
    )

    # Copy the template sheet to the very end of the workbook, then rename it.
    $template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
    $new = $wb.Worksheets.Item($wb.Worksheets.Count)
    $new.Name = $Name

    # Fill in the user-story (ticket ref) cell before the description cell so
    # new shared-string entries are appended in ticket-ref-then-market order,
    # matching how the sheet was originally authored (B2 = "<Country> Market"
    # description, B4 = ticket ref).
    $new.Range("B4").Value = $TicketRef
    $new.Range("B2").Value = $MarketLabel

    # Match the header-row wrapping height used on the other freshly typed
    # sheets.
    $new.Rows.Item(3).RowHeight = 28.8
    $new.Rows.Item(4).RowHeight = 28.8

    # New sheets default to a plain top-left selection spanning the used
    # range instead of the inherited B4 selection.
    $new.Range("A1:D10").Select() | Out-Null

    return $new
}

$russia = Add-MarketSheet "Russia" "NGC-2929/T3297" "Russia Market"
$finland = Add-MarketSheet "Finland" "NGC-3130/T2957" "Finland Market"
$hungary = Add-MarketSheet "Hungary" "NGC-3104/T3006" "Hungary Market"

# The previously active tab (Netherlands) loses its selection; the new last
# sheet (Hungary) becomes the active tab instead.
$hungary.Activate()
$hungary.Range("I16").Select() | Out-Null
